$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 62-75: column A was mistakenly stored as text; make it numeric ---
$ws.Cells.Item(62, 1).Value = 4
$ws.Cells.Item(63, 1).Value = 3
$ws.Cells.Item(64, 1).Value = 3
$ws.Cells.Item(65, 1).Value = 3
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(67, 1).Value = 3
$ws.Cells.Item(68, 1).Value = 3
$ws.Cells.Item(69, 1).Value = 3
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(71, 1).Value = 3
$ws.Cells.Item(72, 1).Value = 3
$ws.Cells.Item(73, 1).Value = 3
$ws.Cells.Item(74, 1).Value = 3
$ws.Cells.Item(75, 1).Value = 3

# --- Append new rows 76-84 parsed from the RSS feed ---

# Row 76
$t = @'
5
'@
$c = $ws.Cells.Item(76, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Experienced Video Editor Needed for Youtube Ad - Upwork
'@
$c = $ws.Cells.Item(76, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Experienced-Video-Editor-Needed-for-Youtube_%7E01336e182e60adbe05?source=rss
'@
$c = $ws.Cells.Item(76, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are seeking a highly skilled video editor to develop a compelling Youtube ad. The ideal candidate will have a strong portfolio demonstrating their ability to create engaging and high-quality video content. As the video editor, you will be responsible for editing raw footage, adding visual effects and transitions, and ensuring the final product meets our brand guidelines. Attention to detail and the ability to work under tight deadlines are a must. If you have a passion for storytelling through video and the skills to bring our ad to life, we would love to hear from you.
Skills required:
- Proficiency in video editing software (e.g., Adobe Premiere Pro, Final Cut Pro)
- Strong understanding of video editing techniques
- Knowledge of audio editing and color correction
- Ability to work collaboratively and take creative direction
Budget
: $20
Posted On
: June 15, 2024 00:09 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects            
Country
: Australia
click to apply

'@
$c = $ws.Cells.Item(76, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are seeking a highly skilled video editor to develop a compelling Youtube ad. The ideal candidate will have a strong portfolio demonstrating their ability to create engaging and high-quality video content. As the video editor, you will be responsible for editing raw footage, adding visual effects and transitions, and ensuring the final product meets our brand guidelines. Attention to detail and the ability to work under tight deadlines are a must. If you have a passion for storytelling through video and the skills to bring our ad to life, we would love to hear from you.<br /><br />
Skills required:<br />
- Proficiency in video editing software (e.g., Adobe Premiere Pro, Final Cut Pro)<br />
- Strong understanding of video editing techniques<br />
- Knowledge of audio editing and color correction<br />
- Ability to work collaboratively and take creative direction<br /><br /><b>Budget</b>: $20
<br /><b>Posted On</b>: June 15, 2024 00:09 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects            <br /><b>Country</b>: Australia
<br /><a href="https://www.upwork.com/jobs/Experienced-Video-Editor-Needed-for-Youtube_%7E01336e182e60adbe05?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(76, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Sat, 15 Jun 2024 00:09:36 +0000
'@
$c = $ws.Cells.Item(76, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Experienced-Video-Editor-Needed-for-Youtube_%7E01336e182e60adbe05?source=rss
'@
$c = $ws.Cells.Item(76, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$20
'@
$c = $ws.Cells.Item(76, 9)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 15, 2024 00:09 UTC
'@
$c = $ws.Cells.Item(76, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(76, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects
'@
$c = $ws.Cells.Item(76, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Australia
'@
$c = $ws.Cells.Item(76, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(76).AutoFit()

# Row 77
$t = @'
5
'@
$c = $ws.Cells.Item(77, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Content Strategy Expert for YouTube Channel Revamp and Optimization - Upwork
'@
$c = $ws.Cells.Item(77, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Content-Strategy-Expert-for-YouTube-Channel-Revamp-and-Optimization_%7E0150f319bb79913d80?source=rss
'@
$c = $ws.Cells.Item(77, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Hello! I'm looking for an expert YouTube strategist to help me maximize the potential of my existing YouTube content and prepare for the launch of my new podcast. I have several podcast interviews that need much more visibility and engagement. My goal is to build momentum using my previous content and seamlessly merge my personal branding with the new podcast.
Responsibilities:
- Develop a strategy to increase visibility and SEO for my YouTube content.
- Advise on whether to maintain the existing YouTube channel under my name or rebrand it to align with my new podcast.
- Implement SEO best practices, including updating titles, descriptions, tags, and thumbnails.
- Suggest ways to leverage older content to create momentum for the new podcast launch.
- Provide guidance on integrating my personal brand with the upcoming podcast content.
- Recommend and implement AI tools for content enhancement.
Goals:
- Improve engagement and visibility of existing content.
- Establish a strong YouTube presence ahead of the new podcast launch.
Requirements:
- Expertise in YouTube content strategy, SEO, and optimization.
- Strong video editing skills with attention to detail.
- Creative skills in graphic design and branding.
- Familiarity with AI tools for content enhancement.
- Ability to provide strategic advice and implement effective solutions.
If you're interested in this project, please provide an overview of your experience and why you're a good fit, examples of previous work, especially with podcasts or interview content, your proposed approach to improving and repurposing existing content, your availability and estimated timeline, and any questions or additional information needed to get started. Thank you for considering this project. 
I look forward to working with you to bring new life to my content and successfully launch my new podcast.
Posted On
: June 15, 2024 00:07 UTC
Category
: Content Strategy
Skills
:Social Media Marketing,     Content Strategy,     YouTube Marketing,     Search Engine Optimization    
Skills
:        Social Media Marketing,                     Content Strategy,                     YouTube Marketing,                     Search Engine Optimization            
Country
: United States
click to apply

'@
$c = $ws.Cells.Item(77, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Hello! I&#039;m looking for an expert YouTube strategist to help me maximize the potential of my existing YouTube content and prepare for the launch of my new podcast. I have several podcast interviews that need much more visibility and engagement. My goal is to build momentum using my previous content and seamlessly merge my personal branding with the new podcast.<br /><br />
Responsibilities:<br /><br />
- Develop a strategy to increase visibility and SEO for my YouTube content.<br />
- Advise on whether to maintain the existing YouTube channel under my name or rebrand it to align with my new podcast.<br />
- Implement SEO best practices, including updating titles, descriptions, tags, and thumbnails.<br />
- Suggest ways to leverage older content to create momentum for the new podcast launch.<br />
- Provide guidance on integrating my personal brand with the upcoming podcast content.<br />
- Recommend and implement AI tools for content enhancement.<br /><br />
Goals:<br /><br />
- Improve engagement and visibility of existing content.<br />
- Establish a strong YouTube presence ahead of the new podcast launch.<br /><br />
Requirements:<br /><br />
- Expertise in YouTube content strategy, SEO, and optimization.<br />
- Strong video editing skills with attention to detail.<br />
- Creative skills in graphic design and branding.<br />
- Familiarity with AI tools for content enhancement.<br />
- Ability to provide strategic advice and implement effective solutions.<br /><br />
If you&#039;re interested in this project, please provide an overview of your experience and why you&#039;re a good fit, examples of previous work, especially with podcasts or interview content, your proposed approach to improving and repurposing existing content, your availability and estimated timeline, and any questions or additional information needed to get started. Thank you for considering this project. <br /><br />
I look forward to working with you to bring new life to my content and successfully launch my new podcast.<br /><br /><br /><b>Posted On</b>: June 15, 2024 00:07 UTC<br /><b>Category</b>: Content Strategy<br /><b>Skills</b>:Social Media Marketing,     Content Strategy,     YouTube Marketing,     Search Engine Optimization    
<br /><b>Skills</b>:        Social Media Marketing,                     Content Strategy,                     YouTube Marketing,                     Search Engine Optimization            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Content-Strategy-Expert-for-YouTube-Channel-Revamp-and-Optimization_%7E0150f319bb79913d80?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(77, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Sat, 15 Jun 2024 00:07:03 +0000
'@
$c = $ws.Cells.Item(77, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Content-Strategy-Expert-for-YouTube-Channel-Revamp-and-Optimization_%7E0150f319bb79913d80?source=rss
'@
$c = $ws.Cells.Item(77, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 15, 2024 00:07 UTC
'@
$c = $ws.Cells.Item(77, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Content Strategy
'@
$c = $ws.Cells.Item(77, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Social Media Marketing,     Content Strategy,     YouTube Marketing,     Search Engine Optimization
'@
$c = $ws.Cells.Item(77, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
United States
'@
$c = $ws.Cells.Item(77, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(77).AutoFit()

# Row 78
$t = @'
5
'@
$c = $ws.Cells.Item(78, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Instagram Video Editing and Scheduling Assistant - Upwork
'@
$c = $ws.Cells.Item(78, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Instagram-Video-Editing-and-Scheduling-Assistant_%7E01759c0bdaefb5e0cf?source=rss
'@
$c = $ws.Cells.Item(78, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
I am looking for a talented individual to assist me LIVE with editing, posting, and scheduling our Instagram videos. Your primary responsibility will be to edit my videos to ensure they are visually appealing and engaging. Additionally, you will be responsible for scheduling the videos to be posted at optimal times to maximize audience reach. This is a great opportunity for someone who wants to learn and grow in the field of social media management.  I would like to LEARN how to eventually do it myself after a few months form an expert.
  Responsibilities:
  - Weekly meetings with me to help me via google meets
  - Edit videos to enhance visual appeal
  - Post videos on Instagram
  - Schedule videos to be posted at optimal times
  Requirements:
  - Proficiency in video editing software
  - Strong understanding of Instagram and its features
  - Excellent organizational and time management skills
  Size: Small
  Duration: 3 to 6 months
  Expertise: Intermediate
Hourly Range
: $6.00-$15.00
Posted On
: June 15, 2024 00:06 UTC
Category
: Video Editing
Skills
:Social Media Marketing,     Social Media Management,     Instagram,     Video Editing    
Skills
:        Social Media Marketing,                     Social Media Management,                     Instagram,                     Video Editing            
Country
: Canada
click to apply

'@
$c = $ws.Cells.Item(78, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
I am looking for a talented individual to assist me LIVE with editing, posting, and scheduling our Instagram videos. Your primary responsibility will be to edit my videos to ensure they are visually appealing and engaging. Additionally, you will be responsible for scheduling the videos to be posted at optimal times to maximize audience reach. This is a great opportunity for someone who wants to learn and grow in the field of social media management.&nbsp;&nbsp;I would like to LEARN how to eventually do it myself after a few months form an expert.<br /><br />
&nbsp;&nbsp;Responsibilities:<br />
&nbsp;&nbsp;- Weekly meetings with me to help me via google meets<br />
&nbsp;&nbsp;- Edit videos to enhance visual appeal<br />
&nbsp;&nbsp;- Post videos on Instagram<br />
&nbsp;&nbsp;- Schedule videos to be posted at optimal times<br /><br />
&nbsp;&nbsp;Requirements:<br />
&nbsp;&nbsp;- Proficiency in video editing software<br />
&nbsp;&nbsp;- Strong understanding of Instagram and its features<br />
&nbsp;&nbsp;- Excellent organizational and time management skills<br /><br />
&nbsp;&nbsp;Size: Small<br />
&nbsp;&nbsp;Duration: 3 to 6 months<br />
&nbsp;&nbsp;Expertise: Intermediate<br /><br /><b>Hourly Range</b>: $6.00-$15.00
<br /><b>Posted On</b>: June 15, 2024 00:06 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Social Media Marketing,     Social Media Management,     Instagram,     Video Editing    
<br /><b>Skills</b>:        Social Media Marketing,                     Social Media Management,                     Instagram,                     Video Editing            <br /><b>Country</b>: Canada
<br /><a href="https://www.upwork.com/jobs/Instagram-Video-Editing-and-Scheduling-Assistant_%7E01759c0bdaefb5e0cf?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(78, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Sat, 15 Jun 2024 00:06:49 +0000
'@
$c = $ws.Cells.Item(78, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Instagram-Video-Editing-and-Scheduling-Assistant_%7E01759c0bdaefb5e0cf?source=rss
'@
$c = $ws.Cells.Item(78, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$6.00-$15.00
'@
$c = $ws.Cells.Item(78, 8)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 15, 2024 00:06 UTC
'@
$c = $ws.Cells.Item(78, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(78, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Social Media Marketing,     Social Media Management,     Instagram,     Video Editing
'@
$c = $ws.Cells.Item(78, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Canada
'@
$c = $ws.Cells.Item(78, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(78).AutoFit()

# Row 79
$t = @'
5
'@
$c = $ws.Cells.Item(79, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Social Video Editor (Descript, Capsule, etc) - Upwork
'@
$c = $ws.Cells.Item(79, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Social-Video-Editor-Descript-Capsule-etc_%7E01657a2a84a0a3698c?source=rss
'@
$c = $ws.Cells.Item(79, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking for a social video editor to add to the team. This would be for simpler projects. Less than 60s videos. Talking heads. Some would involve screen recordings. They would need text and other simple animations.
Expert knowledge of Descript required, as projects will be sourced here, though can use other software to finalize.
Will need to deliver in 16:9; 1:1; 9:16.
This description is a bit broad because there will be many different scenarios. Looking for someone I can trust with all of them, so there is definitely opportunity for long-term work. Price is negotiable once I have more specific details.
Send me your best work.
Budget
: $75
Posted On
: June 14, 2024 23:54 UTC
Category
: Video Editing
Skills
:Descript,     Video Editing    
Skills
:        Descript,                     Video Editing            
Country
: United States
click to apply

'@
$c = $ws.Cells.Item(79, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking for a social video editor to add to the team. This would be for simpler projects. Less than 60s videos. Talking heads. Some would involve screen recordings. They would need text and other simple animations.<br /><br />
Expert knowledge of Descript required, as projects will be sourced here, though can use other software to finalize.<br /><br />
Will need to deliver in 16:9; 1:1; 9:16.<br /><br />
This description is a bit broad because there will be many different scenarios. Looking for someone I can trust with all of them, so there is definitely opportunity for long-term work. Price is negotiable once I have more specific details.<br /><br />
Send me your best work.<br /><br /><b>Budget</b>: $75
<br /><b>Posted On</b>: June 14, 2024 23:54 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Descript,     Video Editing    
<br /><b>Skills</b>:        Descript,                     Video Editing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Social-Video-Editor-Descript-Capsule-etc_%7E01657a2a84a0a3698c?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(79, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:54:10 +0000
'@
$c = $ws.Cells.Item(79, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Social-Video-Editor-Descript-Capsule-etc_%7E01657a2a84a0a3698c?source=rss
'@
$c = $ws.Cells.Item(79, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$75
'@
$c = $ws.Cells.Item(79, 9)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:54 UTC
'@
$c = $ws.Cells.Item(79, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(79, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Descript,     Video Editing
'@
$c = $ws.Cells.Item(79, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
United States
'@
$c = $ws.Cells.Item(79, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(79).AutoFit()

# Row 80
$t = @'
5
'@
$c = $ws.Cells.Item(80, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editor Needed - Upwork
'@
$c = $ws.Cells.Item(80, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-Needed_%7E016c6e7bc4c4a43c58?source=rss
'@
$c = $ws.Cells.Item(80, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are looking for a talented video editor fluent in Hindi, knowledgeable about current trends, and skilled in software like Adobe Premiere Pro or Final Cut Pro. Join our dynamic team to produce high-quality, engaging videos that resonate with our audience.
 Review this https://youtu.be/HBcHlWCDOT0?si=wLPusym-Zoo6I64T for the type of content we need. If you're creative, detail-oriented, and passionate about storytelling through visual media, apply now with your resume and portfolio. We look forward to seeing your work!
Budget
: $50
Posted On
: June 14, 2024 23:51 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects,     Video Production    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects,                     Video Production            
Country
: India
click to apply

'@
$c = $ws.Cells.Item(80, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are looking for a talented video editor fluent in Hindi, knowledgeable about current trends, and skilled in software like Adobe Premiere Pro or Final Cut Pro. Join our dynamic team to produce high-quality, engaging videos that resonate with our audience.<br /><br />
 Review this https://youtu.be/HBcHlWCDOT0?si=wLPusym-Zoo6I64T for the type of content we need. If you&#039;re creative, detail-oriented, and passionate about storytelling through visual media, apply now with your resume and portfolio. We look forward to seeing your work!<br /><br /><b>Budget</b>: $50
<br /><b>Posted On</b>: June 14, 2024 23:51 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects,     Video Production    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects,                     Video Production            <br /><b>Country</b>: India
<br /><a href="https://www.upwork.com/jobs/Video-Editor-Needed_%7E016c6e7bc4c4a43c58?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(80, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:51:41 +0000
'@
$c = $ws.Cells.Item(80, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-Needed_%7E016c6e7bc4c4a43c58?source=rss
'@
$c = $ws.Cells.Item(80, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$50
'@
$c = $ws.Cells.Item(80, 9)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:51 UTC
'@
$c = $ws.Cells.Item(80, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(80, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects,     Video Production
'@
$c = $ws.Cells.Item(80, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
India
'@
$c = $ws.Cells.Item(80, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(80).AutoFit()

# Row 81
$t = @'
5
'@
$c = $ws.Cells.Item(81, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Music Video filmed, edited, colored near Washington DC, United States - Upwork
'@
$c = $ws.Cells.Item(81, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Music-Video-filmed-edited-colored-near-Washington-United-States_%7E01de8d0e98e7748b83?source=rss
'@
$c = $ws.Cells.Item(81, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Multiple Music videos needed - lighting, filmed, edited, colored, synched audio/text
Must be in or able to get to Washington DC
I AM NOT SEEKING SOMEONE WHO RECORDS LIVE PERFORMANCES. 
MUST understand fundamentals of lighting, angles, sequencing/story boarding etc. 
COULD LEAD TO NUMEROUS OTHER VIDEOS
Budget
: $750
Posted On
: June 14, 2024 23:47 UTC
Category
: Videography
Skills
:Drone Videography,     Videography,     Video Production,     Video Post-Editing,     Adobe After Effects,     Music Video,     Video Editing,     Adobe Premiere Pro,     Audio Editing    
Skills
:        Drone Videography,                     Videography,                     Video Production,                     Video Post-Editing,                     Adobe After Effects,                     Music Video,                     Video Editing,                     Adobe Premiere Pro,                     Audio Editing            
Country
: United States
click to apply

'@
$c = $ws.Cells.Item(81, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Multiple Music videos needed - lighting, filmed, edited, colored, synched audio/text<br /><br />
Must be in or able to get to Washington DC<br /><br />
I AM NOT SEEKING SOMEONE WHO RECORDS LIVE PERFORMANCES. <br /><br />
MUST understand fundamentals of lighting, angles, sequencing/story boarding etc. <br /><br />
COULD LEAD TO NUMEROUS OTHER VIDEOS<br /><br /><b>Budget</b>: $750
<br /><b>Posted On</b>: June 14, 2024 23:47 UTC<br /><b>Category</b>: Videography<br /><b>Skills</b>:Drone Videography,     Videography,     Video Production,     Video Post-Editing,     Adobe After Effects,     Music Video,     Video Editing,     Adobe Premiere Pro,     Audio Editing    
<br /><b>Skills</b>:        Drone Videography,                     Videography,                     Video Production,                     Video Post-Editing,                     Adobe After Effects,                     Music Video,                     Video Editing,                     Adobe Premiere Pro,                     Audio Editing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Music-Video-filmed-edited-colored-near-Washington-United-States_%7E01de8d0e98e7748b83?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(81, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:47:31 +0000
'@
$c = $ws.Cells.Item(81, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Music-Video-filmed-edited-colored-near-Washington-United-States_%7E01de8d0e98e7748b83?source=rss
'@
$c = $ws.Cells.Item(81, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$750
'@
$c = $ws.Cells.Item(81, 9)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:47 UTC
'@
$c = $ws.Cells.Item(81, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Videography
'@
$c = $ws.Cells.Item(81, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Drone Videography,     Videography,     Video Production,     Video Post-Editing,     Adobe After Effects,     Music Video,     Video Editing,     Adobe Premiere Pro,     Audio Editing
'@
$c = $ws.Cells.Item(81, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
United States
'@
$c = $ws.Cells.Item(81, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(81).AutoFit()

# Row 82
$t = @'
5
'@
$c = $ws.Cells.Item(82, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editor For Social Media Content - Upwork
'@
$c = $ws.Cells.Item(82, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-For-Social-Media-Content_%7E0149f489e407ad0b46?source=rss
'@
$c = $ws.Cells.Item(82, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking for someone who can edit daily content for social media. Flat rate weekly. Must be creative eye catching visuals.
Budget
: $25
Posted On
: June 14, 2024 23:37 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects            
Country
: United States
click to apply

'@
$c = $ws.Cells.Item(82, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking for someone who can edit daily content for social media. Flat rate weekly. Must be creative eye catching visuals.<br /><br /><b>Budget</b>: $25
<br /><b>Posted On</b>: June 14, 2024 23:37 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Adobe After Effects            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Video-Editor-For-Social-Media-Content_%7E0149f489e407ad0b46?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(82, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:37:52 +0000
'@
$c = $ws.Cells.Item(82, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-For-Social-Media-Content_%7E0149f489e407ad0b46?source=rss
'@
$c = $ws.Cells.Item(82, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$25
'@
$c = $ws.Cells.Item(82, 9)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:37 UTC
'@
$c = $ws.Cells.Item(82, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(82, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Adobe After Effects
'@
$c = $ws.Cells.Item(82, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
United States
'@
$c = $ws.Cells.Item(82, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(82).AutoFit()

# Row 83
$t = @'
5
'@
$c = $ws.Cells.Item(83, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editor for Youtube Channel (in french) - Upwork
'@
$c = $ws.Cells.Item(83, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-french_%7E017703e907ec4c6862?source=rss
'@
$c = $ws.Cells.Item(83, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking to hire an experienced video editor for a Mega Buildings YouTube automation Channel (in French). Videos similar to:
https://www.youtube.com/@TheImpossibleBuild/videos
https://youtu.be/heTD0gNXGM4?si=xzMKx84Vjn7j0uBn
If you have experience editing videos in the Mega Buildings YouTube niche, feel free to apply, and let's talk more details! The plan is to publish two videos per week so time management and adherence to deadlines is a must.
Understanding French will be necessary since the script and voice over will be in French.
This is a long term collaboration with good payment structure.
Looking forward to working with you,
Daniela
Hourly Range
: $10.00-$30.00
Posted On
: June 14, 2024 23:19 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Adobe After Effects,     Video Intro & Outro,     YouTube Development,     French,     English    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Video Production,                     Adobe After Effects,                     Video Intro & Outro,                     YouTube Development,                     French,                     English            
Country
: FRA
click to apply

'@
$c = $ws.Cells.Item(83, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Looking to hire an experienced video editor for a Mega Buildings YouTube automation Channel (in French). Videos similar to:<br />
https://www.youtube.com/@TheImpossibleBuild/videos<br />
https://youtu.be/heTD0gNXGM4?si=xzMKx84Vjn7j0uBn<br /><br />
If you have experience editing videos in the Mega Buildings YouTube niche, feel free to apply, and let&#039;s talk more details! The plan is to publish two videos per week so time management and adherence to deadlines is a must.<br /><br />
Understanding French will be necessary since the script and voice over will be in French.<br /><br />
This is a long term collaboration with good payment structure.<br /><br />
Looking forward to working with you,<br />
Daniela<br /><br /><br /><b>Hourly Range</b>: $10.00-$30.00
<br /><b>Posted On</b>: June 14, 2024 23:19 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Adobe After Effects,     Video Intro &amp; Outro,     YouTube Development,     French,     English    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing,                     Video Production,                     Adobe After Effects,                     Video Intro &amp; Outro,                     YouTube Development,                     French,                     English            <br /><b>Country</b>: FRA
<br /><a href="https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-french_%7E017703e907ec4c6862?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(83, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:19:49 +0000
'@
$c = $ws.Cells.Item(83, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-french_%7E017703e907ec4c6862?source=rss
'@
$c = $ws.Cells.Item(83, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
$10.00-$30.00
'@
$c = $ws.Cells.Item(83, 8)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:19 UTC
'@
$c = $ws.Cells.Item(83, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing
'@
$c = $ws.Cells.Item(83, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Video Editing,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Adobe After Effects,     Video Intro &amp; Outro,     YouTube Development,     French,     English
'@
$c = $ws.Cells.Item(83, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
FRA
'@
$c = $ws.Cells.Item(83, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(83).AutoFit()

# Row 84
$t = @'
5
'@
$c = $ws.Cells.Item(84, 1)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Animated Promotional Video Creation - Upwork
'@
$c = $ws.Cells.Item(84, 2)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Animated-Promotional-Video-Creation_%7E015468e5f296c352a1?source=rss
'@
$c = $ws.Cells.Item(84, 3)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are seeking a talented and experienced video animator to create an engaging and visually stunning promotional video for our company - similar in quality to this: https://www.youtube.com/watch?v=IAs399YLWoY. The video will be used to showcase our products and services, and grab the attention of our target audience. The ideal candidate will have a strong portfolio of animated videos, a creative mindset, and a thorough understanding of video production. 
Skills required:
- Proficiency in video editing software (e.g. Adobe After Effects, Final Cut Pro)
- Ability to create visually appealing and captivating animations
- Attention to detail and ability to meet tight deadlines
We are looking for an intermediate-level animator with a proven track record of delivering high-quality animated videos.
Posted On
: June 14, 2024 23:18 UTC
Category
: Cartoons & Comics
Skills
:Animation,     2D Animation,     Video Production,     Motion Graphics,     Video Commercial    
Skills
:        Animation,                     2D Animation,                     Video Production,                     Motion Graphics,                     Video Commercial            
Country
: United States
click to apply

'@
$c = $ws.Cells.Item(84, 4)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
We are seeking a talented and experienced video animator to create an engaging and visually stunning promotional video for our company - similar in quality to this: https://www.youtube.com/watch?v=IAs399YLWoY. The video will be used to showcase our products and services, and grab the attention of our target audience. The ideal candidate will have a strong portfolio of animated videos, a creative mindset, and a thorough understanding of video production. <br /><br />
Skills required:<br />
- Proficiency in video editing software (e.g. Adobe After Effects, Final Cut Pro)<br />
- Ability to create visually appealing and captivating animations<br />
- Attention to detail and ability to meet tight deadlines<br /><br />
We are looking for an intermediate-level animator with a proven track record of delivering high-quality animated videos.<br /><br /><br /><b>Posted On</b>: June 14, 2024 23:18 UTC<br /><b>Category</b>: Cartoons &amp; Comics<br /><b>Skills</b>:Animation,     2D Animation,     Video Production,     Motion Graphics,     Video Commercial    
<br /><b>Skills</b>:        Animation,                     2D Animation,                     Video Production,                     Motion Graphics,                     Video Commercial            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Animated-Promotional-Video-Creation_%7E015468e5f296c352a1?source=rss">click to apply</a>

'@
$c = $ws.Cells.Item(84, 5)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Fri, 14 Jun 2024 23:18:15 +0000
'@
$c = $ws.Cells.Item(84, 6)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
https://www.upwork.com/jobs/Animated-Promotional-Video-Creation_%7E015468e5f296c352a1?source=rss
'@
$c = $ws.Cells.Item(84, 7)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
June 14, 2024 23:18 UTC
'@
$c = $ws.Cells.Item(84, 10)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Cartoons &amp; Comics
'@
$c = $ws.Cells.Item(84, 11)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
Animation,     2D Animation,     Video Production,     Motion Graphics,     Video Commercial
'@
$c = $ws.Cells.Item(84, 12)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$t = @'
United States
'@
$c = $ws.Cells.Item(84, 13)
$c.NumberFormat = "@"
$c.Value = $t
$c.Style = "Normal"
$ws.Rows.Item(84).AutoFit()
